# Electricity Technology Logit Exponent.xlsx
# Commit: "Fix some pre-run and runtime errors (#232)"
#
# The only substantive content change in the target diff is on the "ETLE"
# worksheet: cell B2 (the logit exponent for "all electricity sources")
# changes from -5 to -4.
#
# (The rest of the diff is Excel application/session metadata - fileVersion
# rupBuild, the author's absolute file path, revisionPtr GUIDs, the
# window position/size in bookViews, and the cached cursor <selection> -
# none of which represent a meaningful spreadsheet edit, so we don't try
# to forge those session-identifier values here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETLE")

$ws.Range("B2").Value = -4
